$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2,4) '63.279.59'
$ws.Cells.Item(2,5).Value2 = '  -5.71%  '

Set-TextValue $ws.Cells.Item(3,4) '3.579.65'
$ws.Cells.Item(3,5).Value2 = '  -0.73%  '

$ws.Cells.Item(4,5).Value2 = '  +0.23%  '

Set-TextValue $ws.Cells.Item(5,4) '401.00'
$ws.Cells.Item(5,5).Value2 = '  -3.60%  '

Set-TextValue $ws.Cells.Item(6,4) '131.04'
$ws.Cells.Item(6,5).Value2 = '  +0.68%  '

Set-TextValue $ws.Cells.Item(7,4) '3.576.56'
$ws.Cells.Item(7,5).Value2 = '  -0.64%  '

Set-TextValue $ws.Cells.Item(8,4) '0.610'
$ws.Cells.Item(8,5).Value2 = '  -6.24%  '

Set-TextValue $ws.Cells.Item(9,4) '1.00'
$ws.Cells.Item(9,5).Value2 = '  +0.09%  '

$ws.Cells.Item(10,5).Value2 = '  -9.13%  '

Set-TextValue $ws.Cells.Item(11,4) '0.157'
$ws.Cells.Item(11,5).Value2 = '  -10.80%  '

Set-TextValue $ws.Cells.Item(12,4) '0.0000310'
$ws.Cells.Item(12,5).Value2 = '  -9.36%  '

Set-TextValue $ws.Cells.Item(13,4) '40.81'
$ws.Cells.Item(13,5).Value2 = '  -4.17%  '

Set-TextValue $ws.Cells.Item(14,4) '9.66'
$ws.Cells.Item(14,5).Value2 = '  -2.47%  '

Set-TextValue $ws.Cells.Item(15,4) '4.138.08'
$ws.Cells.Item(15,5).Value2 = '  -0.67%  '

$ws.Cells.Item(16,5).Value2 = '  -1.48%  '

Set-TextValue $ws.Cells.Item(17,4) '3.572.45'
$ws.Cells.Item(17,5).Value2 = '  -1.21%  '

Set-TextValue $ws.Cells.Item(18,4) '19.60'
$ws.Cells.Item(18,5).Value2 = '  -3.67%  '

Set-TextValue $ws.Cells.Item(19,4) '13.03'
$ws.Cells.Item(19,5).Value2 = '  +5.57%  '

$ws.Cells.Item(20,5).Value2 = '  -6.77%  '

Set-TextValue $ws.Cells.Item(21,4) '63.297.96'
$ws.Cells.Item(21,5).Value2 = '  -5.44%  '

Set-TextValue $ws.Cells.Item(22,4) '412.88'
$ws.Cells.Item(22,5).Value2 = '  -7.98%  '

Set-TextValue $ws.Cells.Item(23,4) '14.71'
$ws.Cells.Item(23,5).Value2 = '  +12.41%  '

Set-TextValue $ws.Cells.Item(24,4) '83.97'
$ws.Cells.Item(24,5).Value2 = '  -5.87%  '

Set-TextValue $ws.Cells.Item(25,4) '2.94'
$ws.Cells.Item(25,5).Value2 = '  -7.18%  '

Set-TextValue $ws.Cells.Item(26,4) '35.01'
$ws.Cells.Item(26,5).Value2 = '  -0.61%  '

Set-TextValue $ws.Cells.Item(27,4) '3.14'
$ws.Cells.Item(27,5).Value2 = '  -6.12%  '

Set-TextValue $ws.Cells.Item(28,4) '9.21'
$ws.Cells.Item(28,5).Value2 = '  -7.75%  '

Set-TextValue $ws.Cells.Item(29,4) '5.14'
$ws.Cells.Item(29,5).Value2 = '  +5.46%  '

Set-TextValue $ws.Cells.Item(30,4) '12.45'
$ws.Cells.Item(30,5).Value2 = '  +0.45%  '

$ws.Cells.Item(31,5).Value2 = '  -2.94%  '

$ws.Cells.Item(32,5).Value2 = '  -3.17%  '

Set-TextValue $ws.Cells.Item(33,4) '6.79'
$ws.Cells.Item(33,5).Value2 = '  -7.89%  '

Set-TextValue $ws.Cells.Item(34,4) '0.156'
$ws.Cells.Item(34,5).Value2 = '  -3.09%  '

Set-TextValue $ws.Cells.Item(35,4) '39.87'
$ws.Cells.Item(35,5).Value2 = '  -0.87%  '

Set-TextValue $ws.Cells.Item(36,4) '1.00'
$ws.Cells.Item(36,5).Value2 = '  +0.00%  '

Set-TextValue $ws.Cells.Item(37,4) '55.27'
$ws.Cells.Item(37,5).Value2 = '  -2.75%  '

$ws.Cells.Item(38,5).Value2 = '  -8.04%  '

Set-TextValue $ws.Cells.Item(39,2) 'ThetaToken'
Set-TextValue $ws.Cells.Item(39,3) 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue $ws.Cells.Item(39,4) '2.82'
$ws.Cells.Item(39,5).Value2 = '  +22.49%  '

Set-TextValue $ws.Cells.Item(40,4) '0.993'
$ws.Cells.Item(40,5).Value2 = '  -0.27%  '

Set-TextValue $ws.Cells.Item(41,2) 'Stellar'
Set-TextValue $ws.Cells.Item(41,3) 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Cells.Item(41,4) '0.138'
$ws.Cells.Item(41,5).Value2 = '  -6.52%  '

Set-TextValue $ws.Cells.Item(42,2) 'ApeXProtocol'
Set-TextValue $ws.Cells.Item(42,3) 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Cells.Item(42,4) '3.12'
$ws.Cells.Item(42,5).Value2 = '  +21.18%  '

Set-TextValue $ws.Cells.Item(43,2) 'Monero'
Set-TextValue $ws.Cells.Item(43,3) 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Cells.Item(43,4) '143.58'
$ws.Cells.Item(43,5).Value2 = '  -3.97%  '

Set-TextValue $ws.Cells.Item(44,2) 'LidoDAOToken'
Set-TextValue $ws.Cells.Item(44,3) 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Cells.Item(44,4) '3.23'
$ws.Cells.Item(44,5).Value2 = '  -1.63%  '

Set-TextValue $ws.Cells.Item(45,2) 'PEPE'
Set-TextValue $ws.Cells.Item(45,3) 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Cells.Item(45,4) '0.0₃0621'
$ws.Cells.Item(45,5).Value2 = '  -14.95%  '

Set-TextValue $ws.Cells.Item(46,2) 'NEARProtocol'
Set-TextValue $ws.Cells.Item(46,3) 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Cells.Item(46,4) '4.29'
$ws.Cells.Item(46,5).Value2 = '  -0.80%  '

Set-TextValue $ws.Cells.Item(47,2) 'EnergySwap'
Set-TextValue $ws.Cells.Item(47,3) 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Cells.Item(47,4) '26.06'
$ws.Cells.Item(47,5).Value2 = '  +21.67%  '

$ws.Cells.Item(48,5).Value2 = '  +2.41%  '

Set-TextValue $ws.Cells.Item(49,4) '2.77'
$ws.Cells.Item(49,5).Value2 = '  -8.34%  '

Set-TextValue $ws.Cells.Item(50,4) '2.50'
$ws.Cells.Item(50,5).Value2 = '  -8.73%  '

Set-TextValue $ws.Cells.Item(51,4) '0.285'
$ws.Cells.Item(51,5).Value2 = '  -9.97%  '
